$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the child-sim timestep duration input (4 day timestep -> 28 day timestep)
$ws.Range("C7").Value = 28

# Introduce the "Formula:" label next to the n-scenarios input
$ws.Range("D11").Value = "Formula:"

# Split out the scenario-count formula: C12 now holds the explicit wave-2 formula
# (219 draws / 3 scenarios-per-draw), and the original n-scenarios-child formula
# moves over to D12.
$ws.Range("D12").Formula = "=C11*C4+1"
$ws.Range("C12").Formula = "=219/3"
